$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '66.413.84'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +2.63%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.186.65'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +1.20%  '

$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.03%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '596.02'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +3.36%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '154.38'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  +3.89%  '

$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.11%  '

$ws.Range('B8').NumberFormat = '@'
$ws.Range('B8').Value = 'XRP'
$ws.Range('C8').NumberFormat = '@'
$ws.Range('C8').Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.551'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +5.09%  '

$ws.Range('B9').NumberFormat = '@'
$ws.Range('B9').Value = 'LidoStakedEther'
$ws.Range('C9').NumberFormat = '@'
$ws.Range('C9').Value = 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '3.176.93'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +1.00%  '

$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +1.27%  '

$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -2.39%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.518'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +4.00%  '

$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +3.32%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '39.25'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +5.92%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.706.55'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +1.04%  '

$ws.Range('B16').NumberFormat = '@'
$ws.Range('B16').Value = 'Polkadot'
$ws.Range('C16').NumberFormat = '@'
$ws.Range('C16').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '7.49'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +5.30%  '

$ws.Range('B17').NumberFormat = '@'
$ws.Range('B17').Value = 'WrappedBTC'
$ws.Range('C17').NumberFormat = '@'
$ws.Range('C17').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '66.370.50'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +2.38%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.187.81'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +0.28%  '

$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +0.64%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '516.39'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +2.47%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '15.39'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +3.73%  '

$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +3.68%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '8.07'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +4.91%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '14.95'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -2.36%  '

$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +2.09%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.00'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +0.00%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.29'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +4.54%  '

$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +3.98%  '

$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  +7.57%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.14'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +15.80%  '

$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +4.60%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '28.34'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +3.34%  '

$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +3.14%  '

$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +0.15%  '

$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +1.37%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '507.37'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +6.57%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '54.92'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +0.60%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0900'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +0.95%  '

$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +2.69%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.128'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +11.12%  '

$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +2.86%  '

$ws.Range('B42').NumberFormat = '@'
$ws.Range('B42').Value = 'dogwifhat'
$ws.Range('C42').NumberFormat = '@'
$ws.Range('C42').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.88'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -1.42%  '

$ws.Range('B43').NumberFormat = '@'
$ws.Range('B43').Value = 'TheGraph'
$ws.Range('C43').NumberFormat = '@'
$ws.Range('C43').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.304'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +8.25%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0₃0674'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +16.30%  '

$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +1.02%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.909.57'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -2.73%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '28.64'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +2.12%  '

$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +3.35%  '

$ws.Range('B50').NumberFormat = '@'
$ws.Range('B50').Value = 'CoreDAO'
$ws.Range('C50').NumberFormat = '@'
$ws.Range('C50').Value = 'https://coinranking.com/coin/HFvoXUQh4+coredao-core'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.65'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +9.96%  '

$ws.Range('B51').NumberFormat = '@'
$ws.Range('B51').Value = 'ThetaToken'
$ws.Range('C51').NumberFormat = '@'
$ws.Range('C51').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.34'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +5.32%  '
